$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value  = "S대 공학 박사가 본 수학 & 통계학이 필요한 이유 – 후기"
$ws.Range("E9").Value  = "https://blog.pabii.co.kr/snu-engineer-phd-math-stat-review/#utm_source=rss&utm_medium=rss&utm_campaign=snu-engineer-phd-math-stat-review"

$ws.Range("D12").Value = "“혼자 공부하는 머신러닝+딥러닝”이 엘리스 아카데미에 등록되었습니다!"
$ws.Range("E12").Value = "https://tensorflow.blog/2021/04/07/%ed%98%bc%ec%9e%90-%ea%b3%b5%eb%b6%80%ed%95%98%eb%8a%94-%eb%a8%b8%ec%8b%a0%eb%9f%ac%eb%8b%9d%eb%94%a5%eb%9f%ac%eb%8b%9d%ec%9d%b4-%ec%97%98%eb%a6%ac%ec%8a%a4-%ec%95%84%ec%b9%b4%eb%8d%b0%eb%af%b8/"

$ws.Range("D16").Value = "CAM (Class activation mapping) 정리 [XAI-2]"
$ws.Range("E16").Value = "https://wewinserv.tistory.com/143"

$ws.Range("D32").Value = "Load balancing을 위한 crontab - 젠킨스 스케쥴러 (h * * * * )"
$ws.Range("E32").Value = "https://dodonam.tistory.com/312"

$ws.Range("D39").Value = "Dimensionality Reduction — Can PCA improve the performance of a classification model?"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Dimensionality-Reduction-%E2%80%94-Can-PCA-improve-the-performance-of-a-classification-model-1"
